{"js": "// Add the new \"SESI\u00d3N VI\" block at the end of the document, mirroring the\n// structure already used for the previous sessions: a heading line, a blank\n// line, the date, the duration, and then three body paragraphs justified\n// (\"both\").\n//\n// The new content is inserted right before the document's final (empty)\n// trailing paragraph \u2014 i.e. right after the last existing paragraph of\n// session V \u2014 so it lands at the very end of the body, ahead of the\n// sectPr.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor on the paragraph right before the very last (trailing empty)\n// paragraph of the document; new content is inserted after it, one\n// paragraph at a time, walking forward so the end result reads top to\n// bottom in the same order as the diff.\nlet anchor = paragraphs.items[paragraphs.items.length - 2];\n\nconst newParagraphs = [\n  { text: \"SESI\u00d3N VI\", justify: false },\n  { text: \"\", justify: false },\n  { text: \"09/06/2025\", justify: false },\n  { text: \"8 horas\", justify: false },\n  {\n    text:\n      \"Durante esta sexta jornada hemos mejorado de forma significativa el sistema de carga de modelos 3D del visor Ennde3D. Se ha sustituido el almacenamiento temporal por localStorage \u2014limitado e inestable para archivos grandes\u2014 por un sistema basado en IndexedDB, m\u00e1s robusto y persistente dentro del navegador. Esto permite cargar archivos de mayor tama\u00f1o de forma m\u00e1s fiable y sin errores al cambiar de p\u00e1gina.\",\n    justify: true,\n  },\n  {\n    text:\n      \"Tambi\u00e9n se ha reorganizado y limpiado el c\u00f3digo anterior, eliminando m\u00e9todos obsoletos como la carga por base64, y se ha unificado el comportamiento de los formatos .glb, .gltf y .stl, incluyendo una rotaci\u00f3n autom\u00e1tica que corrige la orientaci\u00f3n incorrecta de los STL.\",\n    justify: true,\n  },\n  {\n    text:\n      \"En paralelo, se ha redise\u00f1ado la estructura del panel lateral del visor para integrar una cabecera com\u00fan fija con el logo y el bot\u00f3n de acceso al nuevo \u201cModo T\u00e9cnico\u201d. Esta funcionalidad permite alternar entre un panel est\u00e1ndar y uno m\u00e1s avanzado sin superposici\u00f3n de elementos y con transici\u00f3n limpia entre modos.\",\n    justify: true,\n  },\n];\n\nfor (const { text, justify } of newParagraphs) {\n  const p = anchor.insertParagraph(text, Word.InsertLocation.after);\n  // Explicitly (re)set the alignment on every inserted paragraph: Word\n  // otherwise copies the anchor paragraph's formatting (including any\n  // \"both\" justification), which would incorrectly leak onto the plain\n  // heading/date/duration lines.\n  p.alignment = justify ? Word.Alignment.justified : Word.Alignment.left;\n  anchor = p;\n}\n\nawait context.sync();\n", "ps1": "# Add the new \"SESI\u00d3N VI\" block at the end of the document, mirroring the\n# structure already used for the previous sessions: a heading line, a blank\n# line, the date, the duration, and then three body paragraphs justified\n# (\"both\").\n#\n# The new content is inserted right after the last existing paragraph of\n# session V (the empty, \"both\"-justified paragraph that currently sits just\n# before the document's final trailing empty paragraph), so it ends up at\n# the very end of the body, ahead of the sectPr.\n\n$d = $word.ActiveDocument\n\n$wdJustify = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n$wdLeft = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphLeft\n\n# Paragraph right before the document's final (empty) trailing paragraph \u2014\n# the anchor we insert new paragraphs after.\n$anchorIndex = $d.Paragraphs.Count - 1\n$anchor = $d.Paragraphs($anchorIndex)\n\n$durante = \"Durante esta sexta jornada hemos mejorado de forma significativa el sistema de carga de modelos 3D del visor Ennde3D. Se ha sustituido el almacenamiento temporal por localStorage \u2014limitado e inestable para archivos grandes\u2014 por un sistema basado en IndexedDB, m\u00e1s robusto y persistente dentro del navegador. Esto permite cargar archivos de mayor tama\u00f1o de forma m\u00e1s fiable y sin errores al cambiar de p\u00e1gina.\"\n$tambien = \"Tambi\u00e9n se ha reorganizado y limpiado el c\u00f3digo anterior, eliminando m\u00e9todos obsoletos como la carga por base64, y se ha unificado el comportamiento de los formatos .glb, .gltf y .stl, incluyendo una rotaci\u00f3n autom\u00e1tica que corrige la orientaci\u00f3n incorrecta de los STL.\"\n$paralelo = \"En paralelo, se ha redise\u00f1ado la estructura del panel lateral del visor para integrar una cabecera com\u00fan fija con el logo y el bot\u00f3n de acceso al nuevo \u201cModo T\u00e9cnico\u201d. Esta funcionalidad permite alternar entre un panel est\u00e1ndar y uno m\u00e1s avanzado sin superposici\u00f3n de elementos y con transici\u00f3n limpia entre modos.\"\n\n$newParagraphs = @(\n    @{ Text = \"SESI\u00d3N VI\"; Justify = $false },\n    @{ Text = \"\"; Justify = $false },\n    @{ Text = \"09/06/2025\"; Justify = $false },\n    @{ Text = \"8 horas\"; Justify = $false },\n    @{ Text = $durante; Justify = $true },\n    @{ Text = $tambien; Justify = $true },\n    @{ Text = $paralelo; Justify = $true }\n)\n\nforeach ($item in $newParagraphs) {\n    $anchor.Range.InsertParagraphAfter()\n    $anchorIndex = $anchorIndex + 1\n    $newPara = $d.Paragraphs($anchorIndex)\n    if ($item.Text -ne \"\") {\n        $newPara.Range.Text = $item.Text\n    }\n    if ($item.Justify) {\n        $newPara.Alignment = $wdJustify\n    } else {\n        $newPara.Alignment = $wdLeft\n    }\n    $anchor = $newPara\n}\n"}
